$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 60: replace "Number of Data Controller users" / 20 / "Licences"
# with "DC Licence / Support Cost" / 12000 / "Euro"
$ws.Range("C60").Value2 = "DC Licence / Support Cost"
$ws.Range("D60").Value2 = 12000
$ws.Range("E60").Value2 = "Euro"

# Row 62: Total Annual Cost formula now simply mirrors D60 (no more cap/multiplier)
$ws.Range("D62").Formula = "=D60"

# Add a dropdown list data validation on D60 with the allowed licence/support costs
$validation = $ws.Range("D60").Validation
$validation.Delete()
$validation.Add(3, 1, 1, '"0,12000,25000"')
$validation.IgnoreBlank = $true
$validation.InCellDropdown = $true
$validation.ShowInput = $true
$validation.ShowError = $true
